$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Text updates (shared-string content changes)
# ---------------------------------------------------------------------------

# "Ready for handoff" -> "Handed back: in sync with en-US"
# (appears on all three sheets)
$wsOverview.Range("E2:F3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C2:C3").Value     = "Handed back: in sync with en-US"
$wsDeDe.Range("C2:C3").Value     = "Handed back: in sync with en-US"

# Latest Handback DateTime: zh-cn keeps the (updated) shared date text,
# de-de gets its own newer timestamp.
$wsZhCn.Range("K2:K3").Value = "2016-10-14 08:13:44"
$wsDeDe.Range("K2:K3").Value = "2016-10-14 08:14:01"

# ---------------------------------------------------------------------------
# 2. Column width changes
# ---------------------------------------------------------------------------

# Overview: "zh-cn" / "de-de" columns widen
$wsOverview.Range("E:F").ColumnWidth = 29.15

# zh-cn / de-de: "Status" column widens
$wsZhCn.Range("C:C").ColumnWidth = 29.15
$wsDeDe.Range("C:C").ColumnWidth = 29.15

# zh-cn / de-de: "Latest Target File" / "Latest Handback File" widen to 40
$wsZhCn.Range("I:J").ColumnWidth = 39.17
$wsDeDe.Range("I:J").ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# 3. Populate "Latest Target File" (I) / "Latest Handback File" (J) and add
#    hyperlinks for row 2 and row 3 on the zh-cn and de-de sheets.
# ---------------------------------------------------------------------------

$mdUrlPrimary   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/df88b83300403f22a04bea523e71d9080b0b4a73/e2e/39ba4795-f4b7-4eeb-8dcb-0cb950582b83.md"
$mdUrlSecondary = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/df88b83300403f22a04bea523e71d9080b0b4a73/e2e/ffff8c082d66-fb41-4673-b616-f792f6121af7.md"
$mdDisplayPrimary   = "39ba4795-f4b7-4eeb-8dcb-0cb950582b83.md"
$mdDisplaySecondary = "ffff8c082d66-fb41-4673-b616-f792f6121af7.md"

function Rebuild-Hyperlinks($ws) {
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrlPrimary,   "", "", $mdDisplayPrimary)
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdUrlPrimary,   "", "", $mdDisplayPrimary)
    $ws.Hyperlinks.Add($ws.Range("A3"), $mdUrlSecondary, "", "", $mdDisplaySecondary)
    $ws.Hyperlinks.Add($ws.Range("I3"), $mdUrlPrimary,   "", "", $mdDisplayPrimary)
}

# zh-cn sheet
$wsZhCn.Cells.Item(2, 10).Value = "39ba4795-f4b7-4eeb-8dcb-0cb950582b83.b4c75ad0158dd35277e68d189c5ffece562619f6.zh-cn.xlf"
$wsZhCn.Cells.Item(3, 10).Value = "39ba4795-f4b7-4eeb-8dcb-0cb950582b83.b4c75ad0158dd35277e68d189c5ffece562619f6.zh-cn.xlf"
Rebuild-Hyperlinks $wsZhCn

# de-de sheet
$wsDeDe.Cells.Item(2, 10).Value = "39ba4795-f4b7-4eeb-8dcb-0cb950582b83.b4c75ad0158dd35277e68d189c5ffece562619f6.de-de.xlf"
$wsDeDe.Cells.Item(3, 10).Value = "39ba4795-f4b7-4eeb-8dcb-0cb950582b83.b4c75ad0158dd35277e68d189c5ffece562619f6.de-de.xlf"
Rebuild-Hyperlinks $wsDeDe
